$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows that no longer appear in the updated sheet. Clearing the
# contents of a fully-populated single-column row makes the row disappear
# from the saved sheetData (instead of shifting the following rows up),
# which matches the target layout where row numbers 2..43 are preserved
# and only rows 1, 11, 15, 30, 31 are gone.
$ws.Rows(1).ClearContents()
$ws.Rows(11).ClearContents()
$ws.Rows(15).ClearContents()
$ws.Rows(30).ClearContents()
$ws.Rows(31).ClearContents()

# Add the redis-cache reference links to B34 (row for "spring boot 2 redis
# integration") as a hyperlink whose cell text holds all four URLs.
$url = "https://www.journaldev.com/18141/spring-boot-redis-cache"
$linkText = "https://www.journaldev.com/18141/spring-boot-redis-cache`nhttps://dzone.com/articles/implementation-of-redis-in-micro-servicespring-boo`nhttps://www.devglan.com/spring-boot/spring-boot-redis-cache`nhttps://www.concretepage.com/spring-boot/spring-boot-redis-cache"

$ws.Range("B34").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("B34"), $url, "", "", $linkText)
$ws.Rows(34).RowHeight = 57.6

# Restore the view so the previously-edited area (around row 34) is on
# screen with B35 selected, similar to the author's saved view state.
$ws.Range("B35").Select()
